$d = $word.ActiveDocument

# -----------------------------------------------------------------------
# 1) The paragraph right after the table (originally only holding the
#    "_GoBack" bookmark) becomes a plain empty paragraph.
# -----------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# -----------------------------------------------------------------------
# 2) After the "Open the web page URL." step, add two new numbered list
#    steps: "Input user-email." and "Input password.". The "_GoBack"
#    bookmark is re-created at the end of the "Input password." run.
# -----------------------------------------------------------------------
$openPara = $d.Paragraphs(14)

$tail = $openPara.Range
$tail.Collapse(0)
$tail.InsertParagraphAfter()
$tail.InsertParagraphAfter()

# Re-resolve the two freshly inserted (still empty) paragraphs.
$newPara1 = $openPara.Next()
$newPara2 = $newPara1.Next()

$xmlEmail = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:spacing w:line="256" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="TH SarabunPSK" w:hAnsi="TH SarabunPSK" w:cs="TH SarabunPSK"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="TH SarabunPSK" w:hAnsi="TH SarabunPSK" w:cs="TH SarabunPSK" w:hint="cs"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>Input user-email.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$xmlPassword = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:spacing w:line="256" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="TH SarabunPSK" w:hAnsi="TH SarabunPSK" w:cs="TH SarabunPSK"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="TH SarabunPSK" w:hAnsi="TH SarabunPSK" w:cs="TH SarabunPSK" w:hint="cs"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>Input password.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$newPara1.Range.InsertXML($xmlEmail)
$newPara2.Range.InsertXML($xmlPassword)

# -----------------------------------------------------------------------
# 3) Numbering: register a second (unused-by-text but Word-generated)
#    list, "numId=2", mirroring the abstract numbering of numId=1 with a
#    fresh restart override on every level.
# -----------------------------------------------------------------------
